# Add a new attendance column O ("2025-11-20") to the gradebook, mark the
# students who were absent that day with the red "flagged" highlight style
# (the same look previously used on column N), and drop that highlight from
# column N now that O is the newest/most-recent date column. Also normalizes
# N30's stored type from a number to text, matching the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture copy-sources (value+style) before we start mutating anything ---
# E2  = text "0", default style (0)
# K4  = text "1", default style (0)
# J2  = text "0.5", default style (0)
# N1  = header style (1), used as the template for the new O1 header
# N5  = text "1", red-flag style (2)  -- still untouched at this point
# N9  = text "1", red-flag style (2)  -- still untouched at this point
# N30 = number 0, default style (0)   -- still untouched at this point

# New header cell O1: "2025-11-20" (same look as the other header cells).
# Use a leading apostrophe so the dashed string is kept as text instead of
# being auto-parsed into a date serial number.
$ws.Range("N1").Copy($ws.Range("O1")) | Out-Null
$ws.Range("O1").Value = "'2025-11-20"

# Column O body values, rows 2-18 (plain attendance 0/1, text, default style)
$ws.Range("E2").Copy($ws.Range("O2")) | Out-Null
$ws.Range("E2").Copy($ws.Range("O3")) | Out-Null
$ws.Range("E2").Copy($ws.Range("O4")) | Out-Null

# Rows flagged absent on 2025-11-20 get the red highlight style (same as N5/N9/N20 used to have)
$ws.Range("N5").Copy($ws.Range("O5")) | Out-Null
$ws.Range("N5").Copy($ws.Range("O6")) | Out-Null

$ws.Range("E2").Copy($ws.Range("O7")) | Out-Null
$ws.Range("E2").Copy($ws.Range("O8")) | Out-Null

$ws.Range("N9").Copy($ws.Range("O9")) | Out-Null

$ws.Range("E2").Copy($ws.Range("O10")) | Out-Null
$ws.Range("E2").Copy($ws.Range("O11")) | Out-Null

$ws.Range("N5").Copy($ws.Range("O12")) | Out-Null
$ws.Range("N5").Copy($ws.Range("O13")) | Out-Null

$ws.Range("E2").Copy($ws.Range("O14")) | Out-Null
$ws.Range("E2").Copy($ws.Range("O15")) | Out-Null

$ws.Range("N5").Copy($ws.Range("O16")) | Out-Null

$ws.Range("E2").Copy($ws.Range("O17")) | Out-Null
$ws.Range("E2").Copy($ws.Range("O18")) | Out-Null

# Row 19 stores its 0 as a genuine number (matches source data quirk)
$ws.Range("N30").Copy($ws.Range("O19")) | Out-Null

# --- Remove the red-flag highlight from column N now that O carries it ---
$ws.Range("K4").Copy($ws.Range("N5")) | Out-Null
$ws.Range("J2").Copy($ws.Range("N6")) | Out-Null
$ws.Range("K4").Copy($ws.Range("N9")) | Out-Null
$ws.Range("K4").Copy($ws.Range("N20")) | Out-Null

# N30: normalize stored type from number 0 to text "0"
$ws.Range("E2").Copy($ws.Range("N30")) | Out-Null

Write-Output "done"
